# Update per-player serve/return statistics on Sheet1.
# Each assignment below corresponds to a single cell value change
# captured in the source diff (columns are keyed by the header row,
# rows correspond to individual players).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Iga Swiatek (row 2)
$ws.Range("B2").Value = 0.655
$ws.Range("J2").Value = 0.501

# Aryna Sabalenka (row 3)
$ws.Range("J3").Value = 0.368
$ws.Range("L3").Value = 0.447

# Coco Gauff (row 4)
$ws.Range("B4").Value = 0.605
$ws.Range("H4").Value = 0.394

# Elena Rybakina (row 5)
$ws.Range("J5").Value = 0.322

# Jessica Pegula (row 6)
$ws.Range("G6").Value = 0.723
$ws.Range("K6").Value = 0.496

# Marketa Vondrousova (row 8)
$ws.Range("D8").Value = 0.487
$ws.Range("L8").Value = 0.471

# Karolina Muchova (row 9)
$ws.Range("G9").Value = 0.744

# Madison Keys (row 13)
$ws.Range("B13").Value = 0.671
$ws.Range("L13").Value = 0.447

# Sara Sorribes Tormo (row 49)
$ws.Range("C49").Value = 0.588

# Cristina Bucsa (row 62)
$ws.Range("J62").Value = 0.5

# Paula Badosa (row 65)
$ws.Range("G65").Value = 0.74
$ws.Range("K65").Value = 0.45

# Nadia Podoroska (row 68)
$ws.Range("G68").Value = 0.627
$ws.Range("K68").Value = 0.44

# Viktoriya Tomova (row 81)
$ws.Range("J81").Value = 0.411

# Viktorija Golubic (row 84)
$ws.Range("J84").Value = 0.4

# Kayla Day (row 88)
$ws.Range("J88").Value = 0.303
$ws.Range("K88").Value = 0.4

# Zhu Oxuanbai (row 89)
$ws.Range("E89").Value = 0.52
$ws.Range("I89").Value = 0.498

# Linda Fruhvirtova (row 93)
$ws.Range("J93").Value = 0.404

# Oceane Dodin (row 97)
$ws.Range("G97").Value = 0.66

# Claire Liu (row 99)
$ws.Range("L99").Value = 0.401

# Kamilla Rakhimova (row 100)
$ws.Range("J100").Value = 0.38
